$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fix typo "Smarwatch:" -> "Smartwatch:" (insert a "t" after "Smar")
# ------------------------------------------------------------------
$rngTitle = $d.Content
$rngTitle.Find.Execute("Smarwatch:")
$titlePos = $rngTitle.Start + 4
$insTitle = $d.Range($titlePos, $titlePos)
$insTitle.InsertAfter("t")

# ------------------------------------------------------------------
# 2. Extend the table row "#2" description with "(Ruhend oder Sport)"
#    and re-create the _GoBack bookmark inside it (empty, between
#    "Sport" and the closing paren), matching where Word last left the
#    caret after typing.
# ------------------------------------------------------------------
$rngRow = $d.Content
$rngRow.Find.Execute("Abfrage, vor einer Messung, nach Aktivitätszustand ")
$rowInsertPos = $rngRow.End

$insRow = $d.Range($rowInsertPos, $rowInsertPos)
$insRow.InsertAfter("(Ruhend oder Sport)")

# ------------------------------------------------------------------
# 3. Move the _GoBack bookmark: remove it from its old location
#    (after the "QtQuick.Controls 1.2" tab) and place a new, empty one
#    between "Sport" and ")" above.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
}

$newBmPos = $rowInsertPos + 18
$newBmRng = $d.Range($newBmPos, $newBmPos)
$d.Bookmarks.Add("_GoBack", $newBmRng)

# ------------------------------------------------------------------
# 4. Footer page-number field cached result: "2" -> "1"
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$pageChar = $ftr.Range.Characters.Item(1)
$pageChar.Text = "1"
